# Replace "KTKT" with "Hha" in C2:C5 and drop the special Arial/8pt style
# those cells had, so they fall back to the sheet's default formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2..5) {
    $cell = $ws.Cells.Item($r, 3)   # column C
    $cell.ClearFormats()
    $cell.Value = "Hha"
}

# Move/save the active selection to G5, matching the final cursor position.
$ws.Range("G5").Select() | Out-Null
